$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.163.59"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.372.34"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.59"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.46"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.34"
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0787"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.23"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "2.740.79"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "2.363.90"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "43.184.90"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.96"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.80"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.54"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  +7.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.35"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.03"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.60"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +8.67%  "
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "129.24"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.31"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.02"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("D43").Value = "1.931.70"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.24"
$ws.Range("E46").Value = "  -9.14%  "
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "2.599.97"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  -2.63%  "
